$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix first names that had a stray trailing space baked into the
#     shared string ("Bruce ", "Tony ", "Clark ", "Barry ", "Steve ",
#     "Selina ", "Natasha ") - retype them clean so the CONCAT-built
#     email addresses don't contain the extra space either.
$ws.Range("A10").Value = "Selina"
$ws.Range("A11").Value = "Bruce"
$ws.Range("A12").Value = "Tony"
$ws.Range("A13").Value = "Clark"
$ws.Range("A14").Value = "Barry"
$ws.Range("A15").Value = "Steve"
$ws.Range("A16").Value = "Natasha"

# --- Swap the modern CONCAT() function for the legacy CONCATENATE()
#     in every row of the mail_id column so it no longer needs the
#     _xlfn. future-function prefix.
for ($r = 2; $r -le 16; $r++) {
    $formula = '=CONCATENATE(LOWER(A' + $r + '),".",LOWER(B' + $r + '),"@",J' + $r + ')'
    $ws.Range("K$r").Formula = $formula
}

# --- Drop the baby-name hyperlinks that used to sit on column A.
[void]$ws.Hyperlinks.Delete()

# --- Column A no longer needs to stay pinned to its "best fit" width.
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# --- Column J (mail domain helper column) is unhidden again.
$ws.Columns.Item(10).Hidden = $false

# --- Move the active selection from F20 to K20.
[void]$ws.Range("K20").Select()
